$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.562.48'
$ws.Range("E2").Value = '  +1.73%  '
$ws.Range("D3").Value = '1.597.16'
$ws.Range("E3").Value = '  +0.98%  '
$ws.Range("E4").Value = '  +0.43%  '
$ws.Range("D5").Value = '''212.14'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.13%  '
$ws.Range("E6").Value = '  -0.58%  '
$ws.Range("E7").Value = '  +0.41%  '
$ws.Range("D8").Value = '''26.88'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +5.53%  '
$ws.Range("D9").Value = '''43.72'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.55%  '
$ws.Range("E10").Value = '  +1.04%  '
$ws.Range("E11").Value = '  +0.76%  '
$ws.Range("E12").Value = '  +0.84%  '
$ws.Range("D13").Value = '1.825.40'
$ws.Range("E13").Value = '  +0.92%  '
$ws.Range("D14").Value = '1.592.99'
$ws.Range("E14").Value = '  +0.13%  '
$ws.Range("D15").Value = '29.565.94'
$ws.Range("E15").Value = '  +1.60%  '
$ws.Range("E16").Value = '  +2.86%  '
$ws.Range("D17").Value = '''3.73'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.95%  '
$ws.Range("D18").Value = '''63.80'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.32%  '
$ws.Range("D19").Value = '''241.63'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.47%  '
$ws.Range("D20").Value = '''7.58'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.05%  '
$ws.Range("E21").Value = '  -0.16%  '
$ws.Range("E22").Value = '  +0.66%  '
$ws.Range("E23").Value = '  -0.70%  '
$ws.Range("E24").Value = '  +0.43%  '
$ws.Range("E25").Value = '  -0.07%  '
$ws.Range("E26").Value = '  +1.22%  '
$ws.Range("D27").Value = '''15.38'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.55%  '
$ws.Range("E28").Value = '  -0.01%  '
$ws.Range("E29").Value = '  +1.12%  '
$ws.Range("E30").Value = '  +0.46%  '
$ws.Range("D31").Value = '''0.0477'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.80%  '
$ws.Range("E32").Value = '  +0.37%  '
$ws.Range("E33").Value = '  +0.15%  '
$ws.Range("E34").Value = '  +3.36%  '
$ws.Range("D35").Value = '1.429.02'
$ws.Range("E35").Value = '  +0.40%  '
$ws.Range("E36").Value = '  +1.98%  '
$ws.Range("E37").Value = '  -1.63%  '
$ws.Range("D38").Value = '''2.87'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.17%  '
$ws.Range("E39").Value = '  +0.55%  '
$ws.Range("D40").Value = '''0.0166'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.49%  '
$ws.Range("D41").Value = '''0.540'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.76%  '
$ws.Range("B42").Value = 'BitcoinSV'
$ws.Range("C42").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D42").Value = '''54.66'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.95%  '
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").Value = '''1.96'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.88%  '
$ws.Range("B44").Value = 'Kaspa'
$ws.Range("C44").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D44").Value = '''0.0492'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +7.00%  '
$ws.Range("E45").Value = '  +1.68%  '
$ws.Range("E46").Value = '  +0.49%  '
$ws.Range("D47").Value = '''0.993'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +17.11%  '
$ws.Range("D48").Value = '''65.51'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.32%  '
$ws.Range("E49").Value = '  +0.12%  '
$ws.Range("D50").Value = '1.737.36'
$ws.Range("E50").Value = '  +0.99%  '
$ws.Range("D51").Value = '''85.93'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.45%  '
